# This script applies a weekly data update:
# two new rows of "Pepino ensalada" price data (date 2022-10-21, serial 44855)
# are inserted at the top of the data block (rows 300-301), pushing the
# existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block; this shifts the
# existing rows 300:374 down to 302:376 and keeps their formatting/styles.
$ws.Rows("300:301").Insert()

# New row 300 - "Primera" quality
$ws.Range("A300").Value = 1
$ws.Range("B300").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C300").Value = "Arica y Parinacota"
$ws.Range("D300").Value = 44855
$ws.Range("E300").Value = 15
$ws.Range("F300").Value = 100112043
$ws.Range("G300").Value = "Pepino ensalada"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 290
$ws.Range("K300").Value = 14000
$ws.Range("L300").Value = 15000
$ws.Range("M300").Value = 14517
$ws.Range("N300").Value = "$/caja 70 unidades"
$ws.Range("O300").Value = "Región de Arica y Parinacota"
$ws.Range("P300").Value = 207
$ws.Range("Q300").Value = 70
$ws.Range("R300").Value = "Hortaliza"

# New row 301 - "Segunda" quality
$ws.Range("A301").Value = 1
$ws.Range("B301").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C301").Value = "Arica y Parinacota"
$ws.Range("D301").Value = 44855
$ws.Range("E301").Value = 15
$ws.Range("F301").Value = 100112043
$ws.Range("G301").Value = "Pepino ensalada"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("I301").Value = "Segunda"
$ws.Range("J301").Value = 160
$ws.Range("K301").Value = 11000
$ws.Range("L301").Value = 12000
$ws.Range("M301").Value = 11500
$ws.Range("N301").Value = "$/caja 100 unidades"
$ws.Range("O301").Value = "Región de Arica y Parinacota"
$ws.Range("P301").Value = 115
$ws.Range("Q301").Value = 100
$ws.Range("R301").Value = "Hortaliza"
